# Handles float input without breaking stuff
#
# Updates the marksheet with the student's actual scoring data and
# collapses the sheet from a 3-subject layout (A:B, D:E, G:H) down to a
# single subject (A:B) -- the other two question blocks are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Score summary block (rows 10-12)
# ---------------------------------------------------------------------
# Row 10 - "No." (counts)
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 8
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 15
$ws.Range("E10").Value = 28

# Row 11 - "Marking" (per-question marks)
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12 - "Total"
$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 32
$ws.Range("C12").Value = -5
$ws.Range("E12").Value = "27/112"

# ---------------------------------------------------------------------
# Student answers for the remaining subject (column A), rows 16-40
# ---------------------------------------------------------------------
$ws.Range("A16").Value = "Option A"
$ws.Range("A18").Value = "Option B"
$ws.Range("A19").Value = "Option C"
$ws.Range("A22").Value = "Option A"
$ws.Range("A25").Value = "Option A"
$ws.Range("A26").Value = "Option C"
$ws.Range("A27").Value = "Option A"
$ws.Range("A28").Value = "Option B"
$ws.Range("A29").Value = "Option D"
$ws.Range("A33").Value = "Option D"
$ws.Range("A35").Value = "Option C"

# A couple of "Correct Ans" (D column) corrections for the 2nd block,
# before that whole block is dropped below.
$ws.Range("D17").Value = "Option B"
$ws.Range("D18").Value = "Option B"

# ---------------------------------------------------------------------
# Drop the 3rd subject block entirely (columns G:H)
# ---------------------------------------------------------------------
$ws.Range("G:H").Delete()

# ---------------------------------------------------------------------
# Clear out the 2nd subject block (columns D:E) for rows 19-40 -- only
# rows 16-18 keep values (D16 stays blank/no-answer, D17/D18 set above,
# E16/E17/E18 keep their original "Correct Ans" values).
# ---------------------------------------------------------------------
$ws.Range("D19:E40").ClearContents()
